$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.005.38"
$ws.Range("E2").Value = "  +3.55%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.721.47"
$ws.Range("E3").Value = "  +2.65%  "

$ws.Range("E4").Value = "  +0.31%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "614.58"
$ws.Range("E5").Value = "  +10.26%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "193.79"
$ws.Range("E6").Value = "  +15.54%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.639"
$ws.Range("E7").Value = "  +4.78%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.56%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.726"
$ws.Range("E9").Value = "  +5.60%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.162"
$ws.Range("E10").Value = "  +3.08%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "60.20"
$ws.Range("E11").Value = "  +20.62%  "

$ws.Range("E12").Value = "  +3.13%  "

$ws.Range("E13").Value = "  +2.28%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.309.72"
$ws.Range("E14").Value = "  +3.00%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.720.55"
$ws.Range("E15").Value = "  +2.24%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.57"
$ws.Range("E16").Value = "  +3.52%  "

$ws.Range("E17").Value = "  +5.58%  "

$ws.Range("E18").Value = "  +1.24%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.98"
$ws.Range("E19").Value = "  +4.10%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.851.27"
$ws.Range("E20").Value = "  +4.23%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "412.88"
$ws.Range("E21").Value = "  +4.92%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.59"
$ws.Range("E22").Value = "  +6.39%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "90.09"
$ws.Range("E23").Value = "  +5.39%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.10"
$ws.Range("E24").Value = "  +5.17%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.11"
$ws.Range("E25").Value = "  +6.51%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.32"
$ws.Range("E26").Value = "  +9.59%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.83"
$ws.Range("E27").Value = "  +5.34%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.03"
$ws.Range("E28").Value = "  +1.60%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.81"
$ws.Range("E29").Value = "  +6.96%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.91"
$ws.Range("E30").Value = "  +3.82%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.87"
$ws.Range("E31").Value = "  +6.24%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.79"
$ws.Range("E32").Value = "  +5.23%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.123"
$ws.Range("E33").Value = "  +9.37%  "

$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "46.23"
$ws.Range("E34").Value = "  +11.15%  "

$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "644.32"
$ws.Range("E35").Value = "  +12.10%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "67.90"
$ws.Range("E36").Value = "  +6.72%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0838"
$ws.Range("E37").Value = "  -3.74%  "

$ws.Range("E38").Value = "  +9.28%  "

$ws.Range("E39").Value = "  +0.26%  "

$ws.Range("E40").Value = "  +0.55%  "

$ws.Range("E41").Value = "  +9.95%  "

$ws.Range("E42").Value = "  +5.97%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0449"
$ws.Range("E43").Value = "  +6.04%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.64"
$ws.Range("E44").Value = "  +7.01%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.911.45"
$ws.Range("E45").Value = "  +9.98%  "

$ws.Range("E46").Value = "  +7.28%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.31"
$ws.Range("E47").Value = "  +4.56%  "

$ws.Range("E48").Value = "  +5.22%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "145.94"
$ws.Range("E49").Value = "  +3.89%  "

$ws.Range("E50").Value = "  +3.01%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.58"
$ws.Range("E51").Value = "  -11.17%  "

